# "Titel bei ZA angepasst"
#
# Every "Implementierungslogik - anderes Wort finden" placeholder title
# (the working title used while the "Zirkulaere Abhaengigkeiten" slide
# was split up into several detail slides) gets renamed to the final
# title "Zirkulaere Abhaengigkeiten". Also the "Zirkulaere
# Abhaengigkeiten" overview slide (slide 9 / "Hindernisse") drops the
# "Pointer Verhalten unklar, ..." bullet, since that topic now has its
# own dedicated slides.

$p = $ppt.ActivePresentation

$oldTitle = "Implementierungslogik – anderes Wort finden"
$newTitle = "Zirkuläre Abhängigkeiten"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldTitle) {
                $tr.Delete()
                $shape.TextFrame.TextRange.InsertAfter($newTitle)
            }
        }
    }
}

# Slide 9 ("Hindernisse") still lists the "Pointer Verhalten unklar,
# wann Call by Value wann Call by Reference" bullet under "Zirkuläre
# Abhängigkeiten" - remove just that paragraph, the bullets around it
# stay untouched.
$slide9 = $p.Slides.Item(9)
for ($j = 1; $j -le $slide9.Shapes.Count; $j++) {
    $shape = $slide9.Shapes.Item($j)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        for ($k = 1; $k -le $tr.Paragraphs().Count; $k++) {
            $para = $tr.Paragraphs($k, 1)
            if ($para.Text -like "Pointer Verhalten unklar*") {
                $para.Delete()
                break
            }
        }
    }
}
